$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.NumberFormat = "@"
$cD.Value = "63.416.31"
$cD.Style = "Normal"
$ws.Range("E2").Value = "  +1.25%  "
$cD = $ws.Range("D3")
$cD.NumberFormat = "@"
$cD.Value = "2.544.04"
$cD.Style = "Normal"
$ws.Range("E3").Value = "  +4.51%  "
$ws.Range("E4").Value = "  +0.04%  "
$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value = "570.38"
$cD.Style = "Normal"
$ws.Range("E5").Value = "  +2.20%  "
$cD = $ws.Range("D6")
$cD.NumberFormat = "@"
$cD.Value = "150.64"
$cD.Style = "Normal"
$ws.Range("E6").Value = "  +8.04%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.00%  "
$cD = $ws.Range("D9")
$cD.NumberFormat = "@"
$cD.Value = "2.543.00"
$cD.Style = "Normal"
$ws.Range("E9").Value = "  +4.53%  "
$ws.Range("E10").Value = "  +1.46%  "
$cD = $ws.Range("D11")
$cD.NumberFormat = "@"
$cD.Value = "5.69"
$cD.Style = "Normal"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").Value = "  +2.36%  "
$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value = "28.15"
$cD.Style = "Normal"
$ws.Range("E14").Value = "  +7.32%  "
$cD = $ws.Range("D15")
$cD.NumberFormat = "@"
$cD.Value = "2.997.44"
$cD.Style = "Normal"
$ws.Range("E15").Value = "  +4.62%  "
$cD = $ws.Range("D16")
$cD.NumberFormat = "@"
$cD.Value = "63.314.33"
$cD.Style = "Normal"
$ws.Range("E16").Value = "  +1.28%  "
$cD = $ws.Range("D17")
$cD.NumberFormat = "@"
$cD.Value = "0.0000144"
$cD.Style = "Normal"
$ws.Range("E17").Value = "  +1.52%  "
$cD = $ws.Range("D18")
$cD.NumberFormat = "@"
$cD.Value = "2.512.73"
$cD.Style = "Normal"
$ws.Range("E18").Value = "  +3.25%  "
$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value = "11.63"
$cD.Style = "Normal"
$ws.Range("E19").Value = "  +3.23%  "
$cD = $ws.Range("D20")
$cD.NumberFormat = "@"
$cD.Value = "340.32"
$cD.Style = "Normal"
$ws.Range("E20").Value = "  -2.00%  "
$cD = $ws.Range("D21")
$cD.NumberFormat = "@"
$cD.Value = "4.33"
$cD.Style = "Normal"
$ws.Range("E21").Value = "  +2.89%  "
$cD = $ws.Range("D22")
$cD.NumberFormat = "@"
$cD.Value = "6.79"
$cD.Style = "Normal"
$ws.Range("E22").Value = "  -0.47%  "
$cD = $ws.Range("D24")
$cD.NumberFormat = "@"
$cD.Value = "65.91"
$cD.Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("E26").Value = "  +14.56%  "
$cD = $ws.Range("D27")
$cD.NumberFormat = "@"
$cD.Value = "1.60"
$cD.Style = "Normal"
$ws.Range("E27").Value = "  +2.53%  "
$ws.Range("E28").Value = "  -0.04%  "
$cD = $ws.Range("D29")
$cD.NumberFormat = "@"
$cD.Value = "8.45"
$cD.Style = "Normal"
$ws.Range("E29").Value = "  +3.66%  "
$cD = $ws.Range("D30")
$cD.NumberFormat = "@"
$cD.Value = "7.18"
$cD.Style = "Normal"
$ws.Range("E30").Value = "  +10.56%  "
$cD = $ws.Range("D31")
$cD.NumberFormat = "@"
$cD.Value = "0.0₃0821"
$cD.Style = "Normal"
$ws.Range("E31").Value = "  +4.52%  "
$ws.Range("E32").Value = "  +2.06%  "
$cD = $ws.Range("D33")
$cD.NumberFormat = "@"
$cD.Value = "177.88"
$cD.Style = "Normal"
$ws.Range("E33").Value = "  +2.84%  "
$cD = $ws.Range("D34")
$cD.NumberFormat = "@"
$cD.Value = "1.58"
$cD.Style = "Normal"
$ws.Range("E34").Value = "  +8.76%  "
$cD = $ws.Range("D35")
$cD.NumberFormat = "@"
$cD.Value = "423.13"
$cD.Style = "Normal"
$ws.Range("E35").Value = "  +10.90%  "
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("E37").Value = "  +2.18%  "
$cD = $ws.Range("D38")
$cD.NumberFormat = "@"
$cD.Value = "4.43"
$cD.Style = "Normal"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("E41").Value = "  +0.01%  "
$cD = $ws.Range("D42")
$cD.NumberFormat = "@"
$cD.Value = "39.54"
$cD.Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$cD = $ws.Range("D43")
$cD.NumberFormat = "@"
$cD.Value = "153.90"
$cD.Style = "Normal"
$ws.Range("E43").Value = "  +6.20%  "
$cD = $ws.Range("D44")
$cD.NumberFormat = "@"
$cD.Value = "3.78"
$cD.Style = "Normal"
$ws.Range("E44").Value = "  +2.65%  "
$cD = $ws.Range("D45")
$cD.NumberFormat = "@"
$cD.Value = "20.81"
$cD.Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$cD = $ws.Range("D46")
$cD.NumberFormat = "@"
$cD.Value = "0.608"
$cD.Style = "Normal"
$ws.Range("E46").Value = "  +2.00%  "
$cD = $ws.Range("D47")
$cD.NumberFormat = "@"
$cD.Value = "0.0964"
$cD.Style = "Normal"
$ws.Range("E47").Value = "  +1.03%  "
$cD = $ws.Range("D48")
$cD.NumberFormat = "@"
$cD.Value = "0.0523"
$cD.Style = "Normal"
$ws.Range("E48").Value = "  +0.41%  "
$cD = $ws.Range("D49")
$cD.NumberFormat = "@"
$cD.Value = "0.0239"
$cD.Style = "Normal"
$ws.Range("E49").Value = "  +6.93%  "
$cD = $ws.Range("D50")
$cD.NumberFormat = "@"
$cD.Value = "18.53"
$cD.Style = "Normal"
$ws.Range("E50").Value = "  +3.82%  "
$cD = $ws.Range("D51")
$cD.NumberFormat = "@"
$cD.Value = "1.80"
$cD.Style = "Normal"
$ws.Range("E51").Value = "  +4.09%  "
